$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Backplate description text (row 14, column C)
$ws.Range("C14").Value = "Backing plate, 0.125” thick Aluminum with 4-40 threaded holes"

# Fill in the new "feedstock" row for Delrin (row 16) - Qty / Name / Description
$ws.Range("A16").Value = "n/a"
$ws.Range("B16").Value = "Delrin"
$ws.Range("C16").Value = "Delrin (acetal resin), black, 0.125” thick"

# Add a new feedstock row for Aluminum sheet stock (row 18, with row 17 left blank).
# Clone A16's cell formatting (style 3, "Qty" column look) onto A18 first ...
$ws.Range("A16").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
# ... then fill in the feedstock values.
$ws.Range("A18").Value = "n/a"
$ws.Range("B18").Value = "Aluminum"
$ws.Range("C18").Value = "6061 Aluminum sheet, 0.125” thick"

# Move the selection the way the author's session ended up
$ws.Range("C17").Select()
